$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for each data row.
# Update it from 2023-09-16 (serial 45185) to 2023-10-05 (serial 45204)
# for every data row (rows 2 through 158).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 158 }

$newSerial = 45204

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value2 = $newSerial
    }
}
